# League-2019-07-26.xlsx : "Tablib Dataset" sheet rework
#
# - reorders the rows (id descending instead of ascending)
# - renames the TBall abbreviation to "TB"
# - fills in a "description" value for the ITB row
# - adds two new columns: maxLateGames / maxGames
# - re-applies the bold header formatting to the two new header cells
# - keeps the id / maxLateGames / maxGames columns as plain text (not numbers),
#   matching how the sheet was hand-edited in Excel

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (A1:D1 already say id/name/abbreviation/description; add E1/F1)
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "maxLateGames"
$ws.Range("F1").Value = "maxGames"

# Give the two new header cells the same (bold) look as the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Data rows - written out in the new (id-descending) order.
# id / maxLateGames / maxGames are typed in as text, just like "TB"/"None",
# so every one of these helper cells gets the same "@" -> value -> Normal
# treatment; that keeps the stored cell type as text without leaving a
# custom number format behind.
# ---------------------------------------------------------------------------
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$rows = @(
    @{ Row = 2; A = "6"; B = "Major";      C = "Maj";  D = $null;   E = "16"; F = "16" },
    @{ Row = 3; A = "5"; B = "Minor";      C = "Min";  D = $null;   E = "6";  F = "16" },
    @{ Row = 4; A = "4"; B = "PeeWee";     C = "PW";   D = $null;   E = "4";  F = "14" },
    @{ Row = 5; A = "3"; B = "CoachPitch"; C = "CP";   D = $null;   E = "4";  F = "8"  },
    @{ Row = 6; A = "2"; B = "TBall";      C = "TB";   D = $null;   E = "2";  F = "14" },
    @{ Row = 7; A = "1"; B = "ITB";        C = "ITB";  D = "None";  E = "0";  F = "8"  }
)

foreach ($r in $rows) {
    $row = $r.Row

    Set-TextValue $ws.Cells.Item($row, 1) $r.A      # id
    $ws.Cells.Item($row, 2).Value = $r.B             # name
    $ws.Cells.Item($row, 3).Value = $r.C             # abbreviation

    if ($r.D) {
        $ws.Cells.Item($row, 4).Value = $r.D         # description
    } else {
        # keep the description cell present-but-blank, same as the source file
        $ws.Cells.Item($row, 4).NumberFormat = "@"
    }

    Set-TextValue $ws.Cells.Item($row, 5) $r.E       # maxLateGames
    Set-TextValue $ws.Cells.Item($row, 6) $r.F       # maxGames
}

# ---------------------------------------------------------------------------
# Selection, to mirror the saved workbook's cursor position.
# ---------------------------------------------------------------------------
$ws.Range("K17").Select()
